$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "工作表1" to "Sheet1"
# (the _xlnm._FilterDatabase defined name reference updates automatically)
$ws.Name = "Sheet1"

# Add the missing value in column E of row 11
$ws.Range("E11").Value = 6

# Move the active selection on the bottom-left (frozen) pane to E11
$ws.Range("E11").Select()
